$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sample/example data row
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "کالا ۱"
$ws.Range("D2").Value = 120
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "توضیحات"
$ws.Range("H2").Value = "عدد"

# Header row: add new "دسته بندی" (Category) header in column I,
# matching the style already used by the other header cells (H1).
$ws.Range("I1").Value = "دسته بندی"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").Interior.Color = 65535
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").VerticalAlignment = -4108

$ws.Range("I2").Value = "سایر "

# Update selection to match the edited workbook's saved cursor position
$ws.Range("I3").Select()
